# The edit described by the diff is a pure re-ordering of existing data
# rows in the "Artfynd" sheet: whole rows of observation data moved to
# different row numbers (no new species/records were introduced). Three
# independent groups of rows were permuted:
#   Group 1 (rows 18-20): 18<-20, 19<-18, 20<-19
#   Group 2 (rows 21-24): 21<-23, 22<-24, 23<-21, 24<-22
#   Group 3 (rows 25-26): 25<-26, 26<-25
#   Group 4 (rows 38-40): 38<-40, 39<-38, 40<-39
#
# Columns Y/AA (dates) and AD/AE/AG (booleans) are identical across every
# row in each group, so they don't need to be touched. Only the columns
# below actually change value for a given row, so we write exactly those
# cells with the value that the source row held before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group 1: row 18 <- old row 20 ---------------------------------------
$ws.Range("A18").Value = 131199027
$ws.Range("B18").Value = 8451
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 106545
$ws.Range("F18").Value = "Mindre märgborre"
$ws.Range("G18").Value = "Tomicus minor"
$ws.Range("H18").Value = "(Hartig, 1834)"
# "Antal" column is text-typed in this sheet, force text so "3" isn't
# auto-coerced to a number.
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "3"
$ws.Range("J18").Value = ""
$ws.Range("M18").Value = "äldre gnagspår"
$ws.Range("Q18").Value = 485480
$ws.Range("R18").Value = 6783139
$ws.Range("AC18").Value = "rikligt på flera träd, minst tre träd"
$ws.Range("AF18").Value = ""

# --- Group 1: row 19 <- old row 18 ---------------------------------------
$ws.Range("A19").Value = 131202209
$ws.Range("B19").Value = 57881
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 100049
$ws.Range("F19").Value = "Spillkråka"
$ws.Range("G19").Value = "Dryocopus martius"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("J19").Value = ""
$ws.Range("M19").Value = "äldre spår"
$ws.Range("Q19").Value = 485422
$ws.Range("R19").Value = 6783087
$ws.Range("AF19").Value = ""

# --- Group 1: row 20 <- old row 19 ---------------------------------------
$ws.Range("A20").Value = 131202220
$ws.Range("I20").Value = ""
$ws.Range("Q20").Value = 485506
$ws.Range("R20").Value = 6783094
$ws.Range("AC20").Value = ""

# --- Group 2: row 21 <- old row 23 ---------------------------------------
$ws.Range("A21").Value = 131202213
$ws.Range("Q21").Value = 485480
$ws.Range("R21").Value = 6783087

# --- Group 2: row 22 <- old row 24 ---------------------------------------
$ws.Range("A22").Value = 131202661
$ws.Range("Q22").Value = 485435
$ws.Range("R22").Value = 6783121

# --- Group 2: row 23 <- old row 21 ---------------------------------------
$ws.Range("A23").Value = 131199050
$ws.Range("Q23").Value = 485504
$ws.Range("R23").Value = 6783172

# --- Group 2: row 24 <- old row 22 ---------------------------------------
$ws.Range("A24").Value = 131199107
$ws.Range("Q24").Value = 485476
$ws.Range("R24").Value = 6783098

# --- Group 3: row 25 <- old row 26 ---------------------------------------
$ws.Range("A25").Value = 131202538
$ws.Range("Q25").Value = 485500
$ws.Range("R25").Value = 6783141

# --- Group 3: row 26 <- old row 25 ---------------------------------------
$ws.Range("A26").Value = 131198972
$ws.Range("Q26").Value = 485466
$ws.Range("R26").Value = 6783153

# --- Group 4: row 38 <- old row 40 ---------------------------------------
$ws.Range("A38").Value = 131202630
$ws.Range("B38").Value = 57884
$ws.Range("E38").Value = 100109
$ws.Range("F38").Value = "Tretåig hackspett"
$ws.Range("G38").Value = "Picoides tridactylus"
$ws.Range("Q38").Value = 485526
$ws.Range("R38").Value = 6783143

# --- Group 4: row 39 <- old row 38 ---------------------------------------
$ws.Range("A39").Value = 131199044
$ws.Range("Q39").Value = 485494
$ws.Range("R39").Value = 6783163

# --- Group 4: row 40 <- old row 39 ---------------------------------------
$ws.Range("A40").Value = 131198860
$ws.Range("B40").Value = 57881
$ws.Range("E40").Value = 100049
$ws.Range("F40").Value = "Spillkråka"
$ws.Range("G40").Value = "Dryocopus martius"
$ws.Range("Q40").Value = 485442
$ws.Range("R40").Value = 6783096
